$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.003078177322033415
$ws.Range("C2").Value = 9.983522426115931
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 2797.565817734744
$ws.Range("G2").Value = 2826.269215727881

$ws.Range("B3").Value = 0.003078177322033415
$ws.Range("C3").Value = 0.002658071450198252
$ws.Range("D3").Value = 0.1496068669990043
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 0.6887290743729346
